$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 22.65656533333333
$ws.Range("H2").Value = 67.969696
$ws.Range("I2").Value = 0.9268638682343595
$ws.Range("J2").Value = 0.9268638682343595
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.359524333333333
$ws.Range("N2").Value = 22.078573
$ws.Range("O2").Value = 0.6890485322504919
$ws.Range("P2").Value = 0.6890485322504918
$ws.Range("Q2").Value = 166.7415438804231
$ws.Range("R2").Value = 1500.673894923808
$ws.Range("S2").Value = 0.6386541880028987
$ws.Range("T2").Value = 0.6386541880028986
$ws.Range("G3").Value = 22.65656533333333
$ws.Range("H3").Value = 67.969696
$ws.Range("I3").Value = 0.9268638682343595
$ws.Range("J3").Value = 0.9268638682343595
$ws.Range("O3").Value = 0.2252304747913652
$ws.Range("P3").Value = 0.2252304747913652
$ws.Range("Q3").Value = 54.50309424936177
$ws.Range("R3").Value = 490.5278482442559
$ws.Range("S3").Value = 0.2087579891093862
$ws.Range("T3").Value = 0.2087579891093861
$ws.Range("G4").Value = 22.65656533333333
$ws.Range("H4").Value = 67.969696
$ws.Range("I4").Value = 0.9268638682343595
$ws.Range("J4").Value = 0.9268638682343595
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.9155606666666666
$ws.Range("N4").Value = 2.746682
$ws.Range("O4").Value = 0.08572099295814296
$ws.Range("P4").Value = 0.08572099295814296
$ws.Range("Q4").Value = 20.74346006096355
$ws.Range("R4").Value = 186.691140548672
$ws.Range("S4").Value = 0.07945169112207467
$ws.Range("T4").Value = 0.07945169112207467
$ws.Range("I5").Value = 0.04016781697437198
$ws.Range("J5").Value = 0.04016781697437198
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 7.359524333333333
$ws.Range("N5").Value = 22.078573
$ws.Range("O5").Value = 0.6890485322504919
$ws.Range("P5").Value = 0.6890485322504918
$ws.Range("Q5").Value = 7.226135407966444
$ws.Range("R5").Value = 65.03521867169799
$ws.Range("S5").Value = 0.02767757532989741
$ws.Range("T5").Value = 0.0276775753298974
$ws.Range("I6").Value = 0.04016781697437198
$ws.Range("J6").Value = 0.04016781697437198
$ws.Range("O6").Value = 0.2252304747913652
$ws.Range("P6").Value = 0.2252304747913652
$ws.Range("S6").Value = 0.009047016488470462
$ws.Range("T6").Value = 0.00904701648847046
$ws.Range("I7").Value = 0.04016781697437198
$ws.Range("J7").Value = 0.04016781697437198
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.9155606666666666
$ws.Range("N7").Value = 2.746682
$ws.Range("O7").Value = 0.08572099295814296
$ws.Range("P7").Value = 0.08572099295814296
$ws.Range("Q7").Value = 0.8989664347702221
$ws.Range("R7").Value = 8.090697912931999
$ws.Range("S7").Value = 0.003443225156004116
$ws.Range("T7").Value = 0.003443225156004116
$ws.Range("G8").Value = 0.8058883333333333
$ws.Range("H8").Value = 2.417665
$ws.Range("I8").Value = 0.03296831479126849
$ws.Range("J8").Value = 0.03296831479126849
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 7.359524333333333
$ws.Range("N8").Value = 22.078573
$ws.Range("O8").Value = 0.6890485322504919
$ws.Range("P8").Value = 0.6890485322504918
$ws.Range("Q8").Value = 5.930954799116111
$ws.Range("R8").Value = 53.378593192045
$ws.Range("S8").Value = 0.02271676891769574
$ws.Range("T8").Value = 0.02271676891769573
$ws.Range("G9").Value = 0.8058883333333333
$ws.Range("H9").Value = 2.417665
$ws.Range("I9").Value = 0.03296831479126849
$ws.Range("J9").Value = 0.03296831479126849
$ws.Range("O9").Value = 0.2252304747913652
$ws.Range("P9").Value = 0.2252304747913652
$ws.Range("Q9").Value = 1.938661361062778
$ws.Range("R9").Value = 17.447952249565
$ws.Range("S9").Value = 0.007425469193508592
$ws.Range("T9").Value = 0.00742546919350859
$ws.Range("G10").Value = 0.8058883333333333
$ws.Range("H10").Value = 2.417665
$ws.Range("I10").Value = 0.03296831479126849
$ws.Range("J10").Value = 0.03296831479126849
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.9155606666666666
$ws.Range("N10").Value = 2.746682
$ws.Range("O10").Value = 0.08572099295814296
$ws.Range("P10").Value = 0.08572099295814296
$ws.Range("Q10").Value = 0.7378396597255554
$ws.Range("R10").Value = 6.64055693753
$ws.Range("S10").Value = 0.002826076680064167
$ws.Range("T10").Value = 0.002826076680064166
